$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: quantity stays, but template column switches to "MS-all"
$ws.Range("C2").Value = "MS-all"

# Row 3 becomes the "no track" variant (quantity 1)
$ws.Range("A3").Value = "multidiv no track (no track)"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "MS-all"
$ws.Range("D3").Value = "10 (Denver)"
$ws.Range("E3").Value = "OH (On Hand Loc)"
$ws.Range("F3").Value = "Pro-1"

# New row 4: "Serial track" variant (quantity 4)
$ws.Range("A4").Value = "multidiv Serial track (Serial track)"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = "MS-all"
$ws.Range("D4").Value = "10 (Denver)"
$ws.Range("E4").Value = "OH (On Hand Loc)"
$ws.Range("F4").Value = "Pro-1"

# New row 5: the original "Lot and serial track" variant moves here (quantity 3)
$ws.Range("A5").Value = "multidiv serial (Lot and serial track)"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "MS-all"
$ws.Range("D5").Value = "10 (Denver)"
$ws.Range("E5").Value = "OH (On Hand Loc)"
$ws.Range("F5").Value = "Pro-1"

$ws.Range("C2:C3").Select()
